$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: "Finish_Point"/"Finish" become "Next_Point"/"Next" ---
$ws.Range("B10").Value2 = "Next_Point"
$ws.Range("D10").Value2 = "Next"

# --- New row 11: End_Point / End (stage-ending warp) ---
$ws.Range("A11").Value2 = 8
$ws.Range("B11").Value2 = "End_Point"
$ws.Range("C11").Value2 = ""
$ws.Range("D11").Value2 = "End"
$ws.Range("E11").Value2 = 0

# Give the new row the same formatting as the row above it (row 10)
$ws.Range("A10:E10").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Extend the conditional formatting range that covered E8:E10 to E8:E11 ---
$fc = $ws.Range("E8:E10").FormatConditions
$cfRule = $fc.Item(1)
$cfRule.ModifyAppliesToRange($ws.Range("E8:E11")) | Out-Null

# --- Update the saved selection / active cell ---
$ws.Range("G10").Select() | Out-Null
